$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1 (headers): the sheet had been generated with the data row's
# values copied into the header row by mistake - B1:G1 need to become the
# real column-name labels, and the row needs the same trailing metadata
# columns (property_category .. index) the other property sheets carry. ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (data): extend with the same trailing metadata columns. ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "2011-12-22" must stay a literal text value, not get reinterpreted as a
# date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-12-22"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = "徐耀昌"
$ws.Range("L2").Value = 921
$ws.Range("M2").Value = "tmpd3a41"
$ws.Range("N2").Value = 41
